# Deploying to gh-pages: add 2020 data column (Q) to the Tourism GDP table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q, row 4 (year header) - copy formatting from P4 then set the value
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

# New column Q, row 5 (data value) - copy formatting from P5 then set the value
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 3.3

$excel.CutCopyMode = $false

# Move the active selection to match the post-edit workbook state
$ws.Range("R4").Select()
